# Add a new medication entry ("삭센다") to the dictionary list on Sheet1.
# The list lives in column A (A1 header + A2:A56 sorted entries) and is
# covered by an AutoFilter. The new word sorts alphabetically right after
# "삐콤정" (row 28) and before "세레브렉스" (row 29), so insert a fresh row
# at position 29, push everything else down, and fill in the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(29).Insert()
$ws.Range("A29").Value = "삭센다"
